$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.133.76"
$ws.Cells.Item(2, 5).Value = "  -2.02%  "
$ws.Cells.Item(3, 4).Value = "1.562.61"
$ws.Cells.Item(3, 5).Value = "  -1.66%  "
$ws.Cells.Item(4, 5).Value = "  -0.04%  "
$ws.Cells.Item(5, 4).Formula = "'206.06"
$ws.Cells.Item(5, 5).Value = "  -0.51%  "
$ws.Cells.Item(6, 5).Value = "  -1.59%  "
$ws.Cells.Item(7, 5).Value = "  -0.03%  "
$ws.Cells.Item(8, 4).Formula = "'22.19"
$ws.Cells.Item(8, 5).Value = "  -0.15%  "
$ws.Cells.Item(9, 5).Value = "  -1.97%  "
$ws.Cells.Item(10, 5).Value = "  -0.08%  "
$ws.Cells.Item(11, 4).Formula = "'0.0860"
$ws.Cells.Item(11, 5).Value = "  -0.88%  "
$ws.Cells.Item(12, 4).Value = "1.783.47"
$ws.Cells.Item(12, 5).Value = "  -1.73%  "
$ws.Cells.Item(13, 4).Value = "1.557.96"
$ws.Cells.Item(13, 5).Value = "  -1.87%  "
$ws.Cells.Item(14, 4).Formula = "'3.77"
$ws.Cells.Item(14, 5).Value = "  -2.31%  "
$ws.Cells.Item(15, 4).Formula = "'0.516"
$ws.Cells.Item(15, 5).Value = "  -2.93%  "
$ws.Cells.Item(16, 4).Formula = "'62.97"
$ws.Cells.Item(16, 5).Value = "  -0.79%  "
$ws.Cells.Item(17, 4).Value = "27.105.34"
$ws.Cells.Item(17, 5).Value = "  -2.12%  "
$ws.Cells.Item(18, 4).Formula = "'211.98"
$ws.Cells.Item(18, 5).Value = "  -3.63%  "
$ws.Cells.Item(19, 5).Value = "  -1.26%  "
$ws.Cells.Item(20, 4).Formula = "'7.21"
$ws.Cells.Item(20, 5).Value = "  -1.45%  "
$ws.Cells.Item(22, 5).Value = "  -0.80%  "
$ws.Cells.Item(23, 4).Formula = "'9.39"
$ws.Cells.Item(23, 5).Value = "  -2.19%  "
$ws.Cells.Item(24, 5).Value = "  +0.08%  "
$ws.Cells.Item(25, 4).Formula = "'151.81"
$ws.Cells.Item(25, 5).Value = "  -1.26%  "
$ws.Cells.Item(26, 5).Value = "  -4.39%  "
$ws.Cells.Item(27, 4).Formula = "'14.86"
$ws.Cells.Item(27, 5).Value = "  -1.77%  "
$ws.Cells.Item(28, 5).Value = "  -0.03%  "
$ws.Cells.Item(29, 5).Value = "  -1.43%  "
$ws.Cells.Item(30, 5).Value = "  -1.38%  "
$ws.Cells.Item(31, 5).Value = "  -0.86%  "
$ws.Cells.Item(32, 5).Value = "  -1.91%  "
$ws.Cells.Item(33, 4).Value = "1.374.37"
$ws.Cells.Item(33, 5).Value = "  +0.20%  "
$ws.Cells.Item(34, 4).Formula = "'2.94"
$ws.Cells.Item(34, 5).Value = "  +0.24%  "
$ws.Cells.Item(35, 5).Value = "  +0.20%  "
$ws.Cells.Item(36, 4).Formula = "'0.943"
$ws.Cells.Item(36, 5).Value = "  -3.80%  "
$ws.Cells.Item(37, 5).Value = "  -1.07%  "
$ws.Cells.Item(38, 5).Value = "  -1.60%  "
$ws.Cells.Item(39, 2).Value = "ARBITRUM"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(39, 4).Formula = "'0.813"
$ws.Cells.Item(39, 5).Value = "  -1.53%  "
$ws.Cells.Item(40, 2).Value = "ImmutableX"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(40, 4).Formula = "'0.517"
$ws.Cells.Item(40, 5).Value = "  -3.78%  "
$ws.Cells.Item(41, 5).Value = "  +0.01%  "
$ws.Cells.Item(42, 4).Formula = "'0.992"
$ws.Cells.Item(42, 5).Value = "  +1.72%  "
$ws.Cells.Item(43, 4).Formula = "'1.77"
$ws.Cells.Item(43, 5).Value = "  +2.47%  "
$ws.Cells.Item(44, 5).Value = "  -0.08%  "
$ws.Cells.Item(45, 4).Formula = "'63.35"
$ws.Cells.Item(45, 5).Value = "  -1.26%  "
$ws.Cells.Item(46, 4).Formula = "'5.21"
$ws.Cells.Item(46, 5).Value = "  -0.69%  "
$ws.Cells.Item(47, 4).Value = "1.696.95"
$ws.Cells.Item(47, 5).Value = "  -1.68%  "
$ws.Cells.Item(48, 4).Formula = "'85.44"
$ws.Cells.Item(48, 5).Value = "  -2.66%  "
$ws.Cells.Item(49, 4).Value = "0.0₇0992"
$ws.Cells.Item(49, 5).Value = "  -1.27%  "
$ws.Cells.Item(50, 5).Value = "  -0.73%  "
$ws.Cells.Item(51, 5).Value = "  +0.07%  "
